# Updates cryptos list (prices / 1h volume %) to the latest scraped values.
# Some "Price" cells look like plain numbers (e.g. "0.999", "228.69"); a
# leading apostrophe forces Excel to store them as text instead of
# auto-converting to a numeric value, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.896.69"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").Value = "2.092.95"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'228.69"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'60.40"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.385"
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").Value = "'0.0839"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "2.402.72"
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("D14").Value = "'22.00"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "'0.797"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").Value = "'5.47"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "2.107.21"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").Value = "38.796.19"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").Value = "'71.55"
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").Value = "'170.94"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").Value = "'9.51"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  +10.00%  "
$ws.Range("E29").Value = "  +14.53%  "
$ws.Range("D30").Value = "'19.20"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "'0.121"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  +5.58%  "
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").Value = "'4.71"
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "'6.47"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").Value = "'3.62"
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'18.22"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0226"
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.541.49"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'101.06"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").Value = "  +3.64%  "
$ws.Range("D46").Value = "'7.68"
$ws.Range("E46").Value = "  +9.25%  "
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "'4.14"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "2.289.90"
$ws.Range("E51").Value = "  +2.36%  "
